# Horarios Línea 141 - actualización de datos (scrape 13:30:15)
# Hoja "LP1912": datos 141 + hoja "LP1912-215": datos 215 + hoja "6203-6173"
$wb = $excel.ActiveWorkbook

# ---------- Hoja LP1912 ----------
$ws1 = $wb.Worksheets.Item("LP1912")

# Encabezado: hora de actualizacion y total de filas
$ws1.Cells.Item(2, 1).Value = "Última actualización: 13:30:15"
$ws1.Cells.Item(3, 1).Value = "Total filas: 281"

# Filas de datos reordenadas / nuevas (Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada)
$data1 = @(
  @(92, "08:21:27", "09:01", "215A_EL PATO", 40, "LP1912"),
  @(93, "08:21:27", "09:01", "23_HERNANDEZ", 40, "LP1912"),
  @(108, "08:21:27", "09:22", "17_ROMERO", 61, "LP1912"),
  @(109, "07:46:15", "09:22", "16_SANTA ANA", 96, "LP1912"),
  @(110, "07:59:05", "09:23", "16_SANTA ANA", 84, "LP1912"),
  @(111, "07:46:15", "09:23", "17_ROMERO", 97, "LP1912"),
  @(112, "08:21:27", "09:23", "11_ETCHEVERRY", 62, "LP1912"),
  @(206, "11:51:05", "12:13", "10_OLMOS", 22, "LP1912"),
  @(207, "11:51:05", "12:13", "84_COLONIA URQUIZA-ESC 49", 22, "LP1912"),
  @(214, "10:28:12", "12:21", "26_HERNANDEZ", 113, "LP1912"),
  @(215, "10:28:12", "12:21", "215A_EL PATO", 113, "LP1912"),
  @(216, "11:51:05", "12:21", "14_ABASTO", 30, "LP1912"),
  @(217, "12:16:51", "12:21", "16_SANTA ANA", 5, "LP1912"),
  @(224, "10:57:58", "12:37", "17_179 Y 38", 100, "LP1912"),
  @(225, "11:51:05", "12:37", "23_HERNANDEZ", 46, "LP1912"),
  @(226, "11:51:05", "12:37", "27_EL RETIRO", 46, "LP1912"),
  @(230, "12:44:21", "12:44", "10_OLMOS", 0, "LP1912"),
  @(231, "12:44:21", "12:44", "16_SANTA ANA", 0, "LP1912"),
  @(251, "11:51:05", "13:21", "26_HERNANDEZ", 90, "LP1912"),
  @(252, "12:44:21", "13:21", "10_OLMOS", 37, "LP1912"),
  @(259, "13:30:15", "13:38", "10_OLMOS", 8, "LP1912"),
  @(260, "13:30:15", "13:39", "14_ABASTO", 9, "LP1912"),
  @(261, "11:51:05", "13:46", "17_ROMERO", 115, "LP1912"),
  @(262, "12:16:51", "13:50", "215A_EL PATO", 94, "LP1912"),
  @(263, "12:59:34", "13:50", "11_ETCHEVERRY", 51, "LP1912"),
  @(264, "13:30:15", "13:51", "215A_EL PATO", 21, "LP1912"),
  @(265, "13:30:15", "13:51", "11_ETCHEVERRY", 21, "LP1912"),
  @(266, "12:16:51", "13:55", "225_GOMEZ", 99, "LP1912"),
  @(267, "12:44:21", "13:56", "225_GOMEZ", 72, "LP1912"),
  @(268, "12:16:51", "14:04", "17_ROMERO", 108, "LP1912"),
  @(269, "12:44:21", "14:05", "23_HERNANDEZ", 81, "LP1912"),
  @(270, "13:30:15", "14:07", "16_SANTA ANA", 37, "LP1912"),
  @(271, "13:30:15", "14:11", "16_P MOR-167 Y 521", 41, "LP1912"),
  @(272, "12:44:21", "14:13", "16_P MOR-167 Y 521", 89, "LP1912"),
  @(273, "12:44:21", "14:16", "27_EL RETIRO", 77, "LP1912"),
  @(274, "12:44:21", "14:17", "27_EL RETIRO", 93, "LP1912"),
  @(275, "12:44:21", "14:20", "215C_EL PATO", 96, "LP1912"),
  @(276, "12:44:21", "14:21", "26_HERNANDEZ", 97, "LP1912"),
  @(277, "13:30:15", "14:40", "16_SANTA ANA", 70, "LP1912"),
  @(278, "12:59:34", "14:44", "14_ABASTO", 105, "LP1912"),
  @(279, "12:59:34", "14:56", "16_P MOR-SANTA ANA", 117, "LP1912"),
  @(280, "13:30:15", "14:57", "16_P MOR-SANTA ANA", 87, "LP1912"),
  @(281, "12:59:34", "14:58", "215B_EL PATO", 119, "LP1912"),
  @(282, "13:30:15", "15:00", "81_EL PELIGRO", 90, "LP1912"),
  @(283, "13:30:15", "15:05", "10_OLMOS", 95, "LP1912"),
  @(284, "13:30:15", "15:14", "11_ETCHEVERRY", 104, "LP1912"),
  @(285, "13:30:15", "15:20", "15_ABASTO", 110, "LP1912"),
  @(286, "13:30:15", "15:20", "26_HERNANDEZ", 110, "LP1912")
)
foreach ($row in $data1) {
  $n = $row[0]
  $ws1.Cells.Item($n, 1).Value = $row[1]
  $ws1.Cells.Item($n, 2).Value = $row[2]
  $ws1.Cells.Item($n, 3).Value = $row[3]
  $ws1.Cells.Item($n, 4).Value = $row[4]
  $ws1.Cells.Item($n, 5).Value = $row[5]
}

# ---------- Hoja LP1912-215 ----------
$ws2 = $wb.Worksheets.Item("LP1912-215")

# Encabezado: hora de actualizacion y total de filas
$ws2.Cells.Item(2, 1).Value = "Última actualización: 13:30:15"
$ws2.Cells.Item(3, 1).Value = "Total filas: 31"

# Filas de datos reordenadas / nuevas (Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada)
$data2 = @(
  @(34, "13:30:15", "13:51", "215A_EL PATO", 21, "LP1912"),
  @(35, "12:44:21", "14:20", "215C_EL PATO", 96, "LP1912"),
  @(36, "12:59:34", "14:58", "215B_EL PATO", 119, "LP1912")
)
foreach ($row in $data2) {
  $n = $row[0]
  $ws2.Cells.Item($n, 1).Value = $row[1]
  $ws2.Cells.Item($n, 2).Value = $row[2]
  $ws2.Cells.Item($n, 3).Value = $row[3]
  $ws2.Cells.Item($n, 4).Value = $row[4]
  $ws2.Cells.Item($n, 5).Value = $row[5]
}

# ---------- Hoja 6203-6173 ----------
# Solo cambia la hora de "Ultima actualizacion"; los datos y el total de filas no varian
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: 13:30:15"
